# Apply the "c_stock -> c_stocks" rename and associated edits described
# in the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the two template sheets for consistency.
# ---------------------------------------------------------------------
$wsStock    = $wb.Worksheets.Item("c_stock")
$wsStockOld = $wb.Worksheets.Item("c_stock_old")
$wsStock.Name    = "c_stocks"
$wsStockOld.Name = "c_stocks_old"

# ---------------------------------------------------------------------
# 2. c_stocks (formerly c_stock): re-group the B2:B43 "concat" formula
#    into a single shared formula, then recreate the existing shared
#    groups (I39:I43, J39:J43, F40:F43) so they keep working together
#    with the new one.
# ---------------------------------------------------------------------
$quote = [char]34

$wsStock.Range("B2:B43").ClearContents()
$wsStock.Range("I39:J43").ClearContents()
$wsStock.Range("F40:F43").ClearContents()

$wsStock.Range("B2:B43").Formula = "=E2&" + $quote + "_" + $quote + "&C2"
$wsStock.Range("I39:I43").Formula = "=(F39*(1-F39)/(G39*G39)-1)*F39"
$wsStock.Range("J39:J43").Formula = "=(F39*(1-F39)/(G39*G39)-1)*(1-F39)"
$wsStock.Range("F40:F43").Formula = "=1-0.3"

# ---------------------------------------------------------------------
# 3. c_stocks_old (formerly c_stock_old): same re-grouping, this time
#    the B column forms two shared-formula groups (B2:B33, B34:B65)
#    plus the existing I58:I61 group.
# ---------------------------------------------------------------------
$wsStockOld.Range("B2:B33").ClearContents()
$wsStockOld.Range("B34:B61").ClearContents()
$wsStockOld.Range("I58:I61").ClearContents()

$wsStockOld.Range("B2:B33").Formula = "=E2&" + $quote + "_" + $quote + "&C2"
$wsStockOld.Range("B34:B61").Formula = "=E34&" + $quote + "_" + $quote + "&C34"
$wsStockOld.Range("I58:I61").Formula = "=1-0.3"

# Row 51 changes: C51 now stores a new shared string "s" instead of the
# "dg_ev_wet_closed" lookup, so B51's formula recalculates to "DW_s".
$wsStockOld.Range("C51").Value = "s"

# ---------------------------------------------------------------------
# 4. Selection / view bookkeeping for each sheet, finishing with
#    c_stocks_old as the active tab (matches the final activeTab/
#    tabSelected state from the diff).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("user_input_details").Range("B6").Select()
$wb.Worksheets.Item("user_inputs").Range("E20").Select()
$wb.Worksheets.Item("time_periods").Range("D7").Select()
$wb.Worksheets.Item("AD_lu_transitions").Range("A1:N1").Select()
$wsStock.Range("E40").Select()
$wsStockOld.Range("F58").Select()
$wsStockOld.Activate()
